# Cardioid and Koch Parametrics Plots
# -----------------------------------
# The underlying data table on the "Ganho" sheet is trimmed: the rows for
# freq=3 (gain=-0.1) and freq=11 (gain=-0.8) are removed, so the remaining
# eight data points collapse to six (plus the header row) and every row
# below the deletions shifts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ganho")

# Remove the freq=3 row (original row 2) - full row delete shifts rows up.
$ws.Rows.Item(2).Delete()

# After the first delete, the freq=11 row (originally row 6) is now row 5.
$ws.Rows.Item(5).Delete()

# Match the saved selection/active-cell state left on the sheet afterwards.
$ws.Range("B7").Select() | Out-Null
